# Remove the R2 / WDS_ST_EXISTING row (row 83) from the hurricane-damage
# percentage table. This shifts every subsequent row up by one (so the old
# row 140, R3-R2 / TRANSMISSION_INTERREGIONAL, is no longer duplicated and
# the sheet's used range shrinks from A1:L140 to A1:L139).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(83).Delete()
